# BSides22 Project Summary.xlsx - "Cost Summary" sheet update
#
# - Insert a new adjustable-cost line ("Learn to Solder Kits", $400) into the
#   Con Badge cost breakdown, just above "Professional Services" (old row 15,
#   now row 16). This shifts every row from 15 downward by one, which also
#   shifts all of the formulas that reference those rows (handled
#   automatically by Rows.Insert, same as Excel would do).
# - Bump a handful of the already-adjustable dollar figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Summary")

# Insert the new row, pushing "Professional Services" (and everything below)
# down by one row.
$ws.Rows.Item(15).Insert()

# The freshly inserted row inherits a "no borders" default style; copy the
# formatting from the row above (row 14, "Assembly") for each of the three
# column groups so the new row matches the sheet's look (same as Excel's
# Insert behavior when the row above/below share formatting).
$ws.Range("B14:C14").Copy()
$ws.Range("B15:C15").PasteSpecial(-4122)
$ws.Range("E14:F14").Copy()
$ws.Range("E15:F15").PasteSpecial(-4122)
$ws.Range("H14:I14").Copy()
$ws.Range("H15:I15").PasteSpecial(-4122)

# New line item: Learn to Solder Kits - $400 (Con Badge column only).
$ws.Range("B15").Value = "Learn to Solder Kits"
$ws.Range("C15").Value = 400

# Adjustable value updates elsewhere in the breakdown.
$ws.Range("C8").Value = 250   # Prototype Cost (Con Badge)
$ws.Range("I8").Value = 100   # Prototype Cost (SAO)
$ws.Range("C13").Value = 600  # Programming Fee (Con Badge)
$ws.Range("C14").Value = 600  # Assembly (Con Badge)

# View state: zoom in a bit and move the selection.
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("E5").Select()
